$d = $word.ActiveDocument
$d.Content.Find.Execute("78×45=3510", $true, $false, $false, $false, $false, $true, 1, $false, "31×38=1178", 2) | Out-Null
$d.Content.Find.Execute("12×39=468", $true, $false, $false, $false, $false, $true, 1, $false, "50×84=4200", 2) | Out-Null
$d.Content.Find.Execute("39×76=2964", $true, $false, $false, $false, $false, $true, 1, $false, "30×62=1860", 2) | Out-Null
$d.Content.Find.Execute("71×65=4615", $true, $false, $false, $false, $false, $true, 1, $false, "45×92=4140", 2) | Out-Null
$d.Content.Find.Execute("96×19=1824", $true, $false, $false, $false, $false, $true, 1, $false, "17×56=952", 2) | Out-Null
$d.Content.Find.Execute("88×98=8624", $true, $false, $false, $false, $false, $true, 1, $false, "93×20=1860", 2) | Out-Null
$d.Content.Find.Execute("27×80=2160", $true, $false, $false, $false, $false, $true, 1, $false, "22×20=440", 2) | Out-Null
$d.Content.Find.Execute("95×34=3230", $true, $false, $false, $false, $false, $true, 1, $false, "62×41=2542", 2) | Out-Null
$d.Content.Find.Execute("20×41=820", $true, $false, $false, $false, $false, $true, 1, $false, "28×67=1876", 2) | Out-Null
$d.Content.Find.Execute("56×79=4424", $true, $false, $false, $false, $false, $true, 1, $false, "57×32=1824", 2) | Out-Null
$d.Content.Find.Execute("15×32=480", $true, $false, $false, $false, $false, $true, 1, $false, "62×96=5952", 2) | Out-Null
$d.Content.Find.Execute("53×98=5194", $true, $false, $false, $false, $false, $true, 1, $false, "81×57=4617", 2) | Out-Null
$d.Content.Find.Execute("43×89=3827", $true, $false, $false, $false, $false, $true, 1, $false, "99×64=6336", 2) | Out-Null
$d.Content.Find.Execute("79×25=1975", $true, $false, $false, $false, $false, $true, 1, $false, "38×12=456", 2) | Out-Null
$d.Content.Find.Execute("66×37=2442", $true, $false, $false, $false, $false, $true, 1, $false, "77×27=2079", 2) | Out-Null
$d.Content.Find.Execute("19×73=1387", $true, $false, $false, $false, $false, $true, 1, $false, "67×84=5628", 2) | Out-Null
$d.Content.Find.Execute("20×13=260", $true, $false, $false, $false, $false, $true, 1, $false, "99×92=9108", 2) | Out-Null
$d.Content.Find.Execute("27×79=2133", $true, $false, $false, $false, $false, $true, 1, $false, "75×37=2775", 2) | Out-Null
$d.Content.Find.Execute("80×82=6560", $true, $false, $false, $false, $false, $true, 1, $false, "40×21=840", 2) | Out-Null
$d.Content.Find.Execute("42×90=3780", $true, $false, $false, $false, $false, $true, 1, $false, "82×24=1968", 2) | Out-Null
$d.Content.Find.Execute("89×71=6319", $true, $false, $false, $false, $false, $true, 1, $false, "71×82=5822", 2) | Out-Null
$d.Content.Find.Execute("52×69=3588", $true, $false, $false, $false, $false, $true, 1, $false, "93×23=2139", 2) | Out-Null
$d.Content.Find.Execute("53×67=3551", $true, $false, $false, $false, $false, $true, 1, $false, "91×80=7280", 2) | Out-Null
$d.Content.Find.Execute("39×37=1443", $true, $false, $false, $false, $false, $true, 1, $false, "20×17=340", 2) | Out-Null
$d.Content.Find.Execute("57×99=5643", $true, $false, $false, $false, $false, $true, 1, $false, "14×44=616", 2) | Out-Null
